# Generate Report for Handoff
#
# Two newly localized source files were picked up by the pipeline:
#   45d95610-408c-4067-abcf-b6ddc38b996c.md
#   92688179-c346-4db8-a965-aec4e9d62e2e.md
# Each of them gets its own row (with "Ready for handoff" / "Include" status
# and a freshly generated handoff .xlf) inserted just above the permanent
# ".localization-config" row (which is "Not to be localized" / "Ignored"
# and therefore always stays last) on every sheet.

$wb = $excel.ActiveWorkbook

function Style-Hyperlink($ws, $addr) {
    # Re-create the look of the workbook's "HyperLink" cell style
    # (blue, underlined Calibri 11) on the given cell.
    $f = $ws.Range($addr).Font
    $f.Underline = 2
    $f.Color = 15570276
    $f.Name = "Calibri"
    $f.Size = 11
}

# ======================================================================
# Sheet "Overview"
# ======================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md") | Out-Null
Style-Hyperlink $ws "A2"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c463897b4bc5670b746f73360af9d2c57a535b03/e2e/17d16921-d803-4efd-bb22-ef1ca06a2a3e.md", "", "", "17d16921-d803-4efd-bb22-ef1ca06a2a3e.md") | Out-Null
Style-Hyperlink $ws "A3"

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/45d95610-408c-4067-abcf-b6ddc38b996c.md", "", "", "45d95610-408c-4067-abcf-b6ddc38b996c.md") | Out-Null
Style-Hyperlink $ws "A4"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/92688179-c346-4db8-a965-aec4e9d62e2e.md", "", "", "92688179-c346-4db8-a965-aec4e9d62e2e.md") | Out-Null
Style-Hyperlink $ws "A5"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/.localization-config", "", "", ".localization-config") | Out-Null
Style-Hyperlink $ws "A6"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("C6").Value = "Not to be localized"

# ======================================================================
# Sheet "zh-cn"
# ======================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md") | Out-Null
Style-Hyperlink $ws "A2"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/26d8541d6230e3f3e5ecaf821cd8448033cfcc1a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.zh-cn.xlf", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws "C2"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d150d9ea1883f38a8d762011d5e4621677981f79/e2e/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md") | Out-Null
Style-Hyperlink $ws "E2"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5efbe5140c076aca04919483165d9f091539b11e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.zh-cn.xlf", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws "F2"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c463897b4bc5670b746f73360af9d2c57a535b03/e2e/17d16921-d803-4efd-bb22-ef1ca06a2a3e.md", "", "", "17d16921-d803-4efd-bb22-ef1ca06a2a3e.md") | Out-Null
Style-Hyperlink $ws "A3"
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8341175a15e72d81211d1263b4bae898989241a1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.zh-cn.xlf", "", "", "17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws "C3"

# Row 4: 45d95610-... source file
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/45d95610-408c-4067-abcf-b6ddc38b996c.md", "", "", "45d95610-408c-4067-abcf-b6ddc38b996c.md") | Out-Null
Style-Hyperlink $ws "A4"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5efbe5140c076aca04919483165d9f091539b11e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/45d95610-408c-4067-abcf-b6ddc38b996c.f15c29609dc7d85ba2a361885d8790ae332730e5.zh-cn.xlf", "", "", "45d95610-408c-4067-abcf-b6ddc38b996c.f15c29609dc7d85ba2a361885d8790ae332730e5.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws "C4"
$ws.Range("D4").Value = "2016-03-10 03:11:52"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

# Row 5 (new row): 92688179-... source file
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/92688179-c346-4db8-a965-aec4e9d62e2e.md", "", "", "92688179-c346-4db8-a965-aec4e9d62e2e.md") | Out-Null
Style-Hyperlink $ws "A5"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5efbe5140c076aca04919483165d9f091539b11e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/92688179-c346-4db8-a965-aec4e9d62e2e.a6311974d0cebe44d79fa7b57bcbcb4dd831080d.zh-cn.xlf", "", "", "92688179-c346-4db8-a965-aec4e9d62e2e.a6311974d0cebe44d79fa7b57bcbcb4dd831080d.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws "C5"
$ws.Range("D5").Value = "2016-03-10 03:11:52"
$ws.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

# Row 6 (new row): ".localization-config", always last
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/.localization-config", "", "", ".localization-config") | Out-Null
Style-Hyperlink $ws "A6"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("D6").Value = "0001-01-01 00:00:00"
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G6").Value = "0001-01-01 00:00:00"
$ws.Range("H6").Value = "Ignored"

# ======================================================================
# Sheet "de-de"
# ======================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md") | Out-Null
Style-Hyperlink $ws "A2"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95256d14e2226de8f26c490d832cdc5efc7bc2e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.de-de.xlf", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.de-de.xlf") | Out-Null
Style-Hyperlink $ws "C2"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0d4e0f06389a93e0913663c0c75dfc3ef3558808/e2e/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md") | Out-Null
Style-Hyperlink $ws "E2"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a3b35e3481360cce2aa1b967116c08e76da7a2da/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.de-de.xlf", "", "", "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.de-de.xlf") | Out-Null
Style-Hyperlink $ws "F2"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c463897b4bc5670b746f73360af9d2c57a535b03/e2e/17d16921-d803-4efd-bb22-ef1ca06a2a3e.md", "", "", "17d16921-d803-4efd-bb22-ef1ca06a2a3e.md") | Out-Null
Style-Hyperlink $ws "A3"
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d94394ebae535fe0b9027c2b4b735700a647dc7a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.de-de.xlf", "", "", "17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.de-de.xlf") | Out-Null
Style-Hyperlink $ws "C3"

# Row 4: 45d95610-... source file
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/45d95610-408c-4067-abcf-b6ddc38b996c.md", "", "", "45d95610-408c-4067-abcf-b6ddc38b996c.md") | Out-Null
Style-Hyperlink $ws "A4"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a3b35e3481360cce2aa1b967116c08e76da7a2da/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/45d95610-408c-4067-abcf-b6ddc38b996c.f15c29609dc7d85ba2a361885d8790ae332730e5.de-de.xlf", "", "", "45d95610-408c-4067-abcf-b6ddc38b996c.f15c29609dc7d85ba2a361885d8790ae332730e5.de-de.xlf") | Out-Null
Style-Hyperlink $ws "C4"
$ws.Range("D4").Value = "2016-03-10 03:11:55"
$ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

# Row 5 (new row): 92688179-... source file
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/e2e/92688179-c346-4db8-a965-aec4e9d62e2e.md", "", "", "92688179-c346-4db8-a965-aec4e9d62e2e.md") | Out-Null
Style-Hyperlink $ws "A5"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a3b35e3481360cce2aa1b967116c08e76da7a2da/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/92688179-c346-4db8-a965-aec4e9d62e2e.a6311974d0cebe44d79fa7b57bcbcb4dd831080d.de-de.xlf", "", "", "92688179-c346-4db8-a965-aec4e9d62e2e.a6311974d0cebe44d79fa7b57bcbcb4dd831080d.de-de.xlf") | Out-Null
Style-Hyperlink $ws "C5"
$ws.Range("D5").Value = "2016-03-10 03:11:55"
$ws.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

# Row 6 (new row): ".localization-config", always last
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/618c0a4645f07130e350ef23f48060c6fe73969c/.localization-config", "", "", ".localization-config") | Out-Null
Style-Hyperlink $ws "A6"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("D6").Value = "0001-01-01 00:00:00"
$ws.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("G6").Value = "0001-01-01 00:00:00"
$ws.Range("H6").Value = "Ignored"

Write-Host "Localization status report regenerated for handoff."
